# Auto-generated edit script applying numeric corrections to the
# Spriggan_Profits leve-profit tables (per scheduled-runner price refresh).
$wb = $excel.ActiveWorkbook

# ================= Sheet: ALC =================
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 158456.73
$ws.Range("I15").Value = 158456.73
$ws.Range("K15").Value = 475370.1900000001
$ws.Range("M15").Value = -475201.1900000001

# Row 19
$ws.Range("H19").Value = 627.16
$ws.Range("J19").Value = 684.6667
$ws.Range("L19").Value = 684.6667
$ws.Range("N19").Value = -1034.6667

# Row 33
$ws.Range("H33").Value = 811.05
$ws.Range("I33").Value = 621.82355
$ws.Range("K33").Value = 621.82355
$ws.Range("M33").Value = -392.82355

# Row 38
$ws.Range("H38").Value = 1806.4286
$ws.Range("I38").Value = 1129
$ws.Range("K38").Value = 3387
$ws.Range("M38").Value = -3015

# Row 40
$ws.Range("H40").Value = 4638248
$ws.Range("I40").Value = 4139.7646
$ws.Range("K40").Value = 4139.7646
$ws.Range("M40").Value = -3964.7646

# Row 58
$ws.Range("H58").Value = 5538.385
$ws.Range("I58").Value = 216.85715
$ws.Range("J58").Value = 11746.833
$ws.Range("K58").Value = 650.5714499999999
$ws.Range("L58").Value = 35240.499
$ws.Range("M58").Value = -500.5714499999999
$ws.Range("N58").Value = -35540.499

# Row 80
$ws.Range("H80").Value = 722.3182
$ws.Range("I80").Value = 578.2143
$ws.Range("J80").Value = 974.5
$ws.Range("K80").Value = 1734.6429
$ws.Range("L80").Value = 2923.5
$ws.Range("M80").Value = -736.6428999999998
$ws.Range("N80").Value = -4919.5

# Row 83
$ws.Range("H83").Value = 722.3182
$ws.Range("I83").Value = 578.2143
$ws.Range("J83").Value = 974.5
$ws.Range("K83").Value = 5203.928699999999
$ws.Range("L83").Value = 8770.5
$ws.Range("M83").Value = -211.9286999999995
$ws.Range("N83").Value = -18754.5

# Row 92
$ws.Range("H92").Value = 619.25
$ws.Range("I92").Value = 766.38464
$ws.Range("K92").Value = 766.38464
$ws.Range("M92").Value = 481.61536

# Row 103
$ws.Range("H103").Value = 1968
$ws.Range("I103").Value = 1769.6666
$ws.Range("K103").Value = 5308.9998
$ws.Range("M103").Value = -4722.9998

# Row 106
$ws.Range("H106").Value = 2891.6667
$ws.Range("I106").Value = 2770.2
$ws.Range("K106").Value = 2770.2
$ws.Range("M106").Value = -2139.2

# Row 111
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").ClearContents()


# ================= Sheet: ARM =================
$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 3170.2856
$ws.Range("J88").Value = 3348.6667
$ws.Range("L88").Value = 3348.6667
$ws.Range("N88").Value = -4160.6667

# Row 91
$ws.Range("H91").Value = 3170.2856
$ws.Range("J91").Value = 3348.6667
$ws.Range("L91").Value = 3348.6667
$ws.Range("N91").Value = -6156.6667

# Row 97
$ws.Range("H97").Value = 461
$ws.Range("I97").Value = 448
$ws.Range("K97").Value = 448
$ws.Range("M97").Value = 48


# ================= Sheet: BSM =================
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 145081.42
$ws.Range("J107").Value = 335996.34
$ws.Range("L107").Value = 335996.34
$ws.Range("N107").Value = -339836.34

# Row 134
$ws.Range("H134").Value = 11366077
$ws.Range("I134").Value = 13515237
$ws.Range("J134").Value = 6228.4287
$ws.Range("K134").Value = 40545711
$ws.Range("L134").Value = 18685.2861
$ws.Range("M134").Value = -40543176
$ws.Range("N134").Value = -23755.2861


# ================= Sheet: CUL =================
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2064832.5
$ws.Range("I4").Value = 1183483.6
$ws.Range("K4").Value = 3550450.8
$ws.Range("M4").Value = -3550338.8

# Row 40
$ws.Range("H40").Value = 72.818184
$ws.Range("J40").Value = 69.5
$ws.Range("L40").Value = 278
$ws.Range("N40").Value = -416

# Row 97
$ws.Range("H97").Value = 1401.1875
$ws.Range("I97").Value = 94.57143000000001
$ws.Range("J97").Value = 2417.4443
$ws.Range("K97").Value = 283.71429
$ws.Range("L97").Value = 7252.3329
$ws.Range("M97").Value = 212.28571
$ws.Range("N97").Value = -8244.332900000001

# Row 117
$ws.Range("H117").Value = 2344.889
$ws.Range("J117").Value = 2344.889
$ws.Range("L117").Value = 7034.667
$ws.Range("N117").Value = -13918.667


# ================= Sheet: GSM =================
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 347.27274
$ws.Range("I2").Value = 378.27777
$ws.Range("K2").Value = 378.27777
$ws.Range("M2").Value = -265.27777

# Row 57
$ws.Range("H57").Value = 56199.8
$ws.Range("J57").Value = 56199.8
$ws.Range("L57").Value = 56199.8
$ws.Range("N57").Value = -57839.8

# Row 126
$ws.Range("H126").Value = 4575.129
$ws.Range("I126").Value = 4575.129
$ws.Range("K126").Value = 13725.387
$ws.Range("M126").Value = -11255.387

# Row 132
$ws.Range("H132").Value = 7354557
$ws.Range("I132").Value = 8334898
$ws.Range("K132").Value = 25004694
$ws.Range("M132").Value = -25002164


# ================= Sheet: LTW =================
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2717.4285
$ws.Range("I16").Value = 738.5454999999999
$ws.Range("J16").Value = 4894.2
$ws.Range("K16").Value = 738.5454999999999
$ws.Range("L16").Value = 4894.2
$ws.Range("M16").Value = -568.5454999999999
$ws.Range("N16").Value = -5234.2

# Row 35
$ws.Range("H35").Value = 11124.625
$ws.Range("I35").Value = 8499.5
$ws.Range("K35").Value = 8499.5
$ws.Range("M35").Value = -8163.5

# Row 39
$ws.Range("H39").Value = 400
$ws.Range("J39").Value = 400
$ws.Range("L39").Value = 400
$ws.Range("N39").Value = -1320

# Row 45
$ws.Range("H45").Value = 22749.25
$ws.Range("I45").Value = 13666
$ws.Range("J45").Value = 49999
$ws.Range("K45").Value = 13666
$ws.Range("L45").Value = 49999
$ws.Range("M45").Value = -13259
$ws.Range("N45").Value = -50813

# Row 46
$ws.Range("H46").Value = 891.1667
$ws.Range("I46").Value = 891.1667
$ws.Range("K46").Value = 891.1667
$ws.Range("M46").Value = -703.1667

# Row 55
$ws.Range("H55").Value = 190.58333
$ws.Range("I55").Value = 129
$ws.Range("K55").Value = 129
$ws.Range("M55").Value = 44

# Row 68
$ws.Range("H68").Value = 987
$ws.Range("I68").Value = 987
$ws.Range("K68").Value = 987
$ws.Range("M68").Value = -238

# Row 71
$ws.Range("H71").Value = 987
$ws.Range("I71").Value = 987
$ws.Range("K71").Value = 4935
$ws.Range("M71").Value = -1191

# Row 136
$ws.Range("H136").Value = 2784.05
$ws.Range("I136").Value = 1938.8
$ws.Range("J136").Value = 3065.8
$ws.Range("K136").Value = 5816.4
$ws.Range("L136").Value = 9197.400000000001
$ws.Range("M136").Value = -3266.4
$ws.Range("N136").Value = -14297.4


# ================= Sheet: WVR =================
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 10857.667
$ws.Range("J45").Value = 12829.2
$ws.Range("L45").Value = 12829.2
$ws.Range("N45").Value = -13811.2

# Row 62
$ws.Range("H62").Value = 13158.5
$ws.Range("I62").Value = 7737.75
$ws.Range("K62").Value = 7737.75
$ws.Range("M62").Value = -7113.75

# Row 65
$ws.Range("H65").Value = 13158.5
$ws.Range("I65").Value = 7737.75
$ws.Range("K65").Value = 38688.75
$ws.Range("M65").Value = -35568.75

# Row 122
$ws.Range("H122").Value = 3206.25
$ws.Range("I122").Value = 3131.9565
$ws.Range("K122").Value = 9395.869499999999
$ws.Range("M122").Value = -6945.869499999999

# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

